$wb = $excel.ActiveWorkbook

# Both the "展览" and "全部类型" worksheets contain the same data table and
# both need the "想去人数" (column F) values updated identically.
$sheetNames = @("展览", "全部类型")

# Mapping of row -> new value for column F
$updates = @{
    2  = 1166
    3  = 110
    4  = 1609
    5  = 619
    8  = 11452
    12 = 357
    13 = 1091
    14 = 795
    15 = 12368
    16 = 13045
    21 = 218
    24 = 107
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
